$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 331, shifting existing rows 331-357 down to 332-358.
$ws.Rows(331).Insert()

# Populate the new row 331 with its data. Columns A,B,C,E,F,G,H,I,N,O,Q,R
# carry the same constant values as the surrounding rows for this subset;
# D,J,K,L,M,P hold the new record's values.
$ws.Cells.Item(331, 1).Value = 8
$ws.Cells.Item(331, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(331, 3).Value = "Coquimbo"
$ws.Cells.Item(331, 4).Value = 44931
$ws.Cells.Item(331, 5).Value = 4
$ws.Cells.Item(331, 6).Value = 100112012
$ws.Cells.Item(331, 7).Value = "Espinaca"
$ws.Cells.Item(331, 8).Value = "Sin especificar"
$ws.Cells.Item(331, 9).Value = "Primera"
$ws.Cells.Item(331, 10).Value = 1360
$ws.Cells.Item(331, 11).Value = 500
$ws.Cells.Item(331, 12).Value = 600
$ws.Cells.Item(331, 13).Value = 550
$ws.Cells.Item(331, 14).Value = "$/atado 300 a 500 gramos"
$ws.Cells.Item(331, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(331, 16).Value = 1100
$ws.Cells.Item(331, 17).Value = 0.5
$ws.Cells.Item(331, 18).Value = "Hortaliza"
